$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.820.74"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.642.68"
$ws.Range("E3").Value = "  +0.32%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.32%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.70"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.34%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +1.02%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0619"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -0.43%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.873.17"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.661.75"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  +0.32%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.17"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "26.864.32"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.54%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.09"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("E20").Value = "  -0.29%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +7.55%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +6.64%  "
$ws.Range("E24").Value = "  -1.25%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.09"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.32%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -0.34%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.32"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "1.241.76"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.99%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +3.06%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "1.786.14"
$ws.Range("E43").Value = "  +0.62%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("E45").Value = "  +1.25%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.62"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  +1.04%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.56"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.43%  "
